$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 332, shifting existing rows 332-446 down to 333-447.
$ws.Rows.Item(332).Insert()

# Populate the newly inserted row 332 with the new weekly data point.
$ws.Cells.Item(332, 1).Value = 7
$ws.Cells.Item(332, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(332, 3).Value = "Ñuble"
$ws.Cells.Item(332, 4).Value = 44900
$ws.Cells.Item(332, 5).Value = 16
$ws.Cells.Item(332, 6).Value = 100114001
$ws.Cells.Item(332, 7).Value = "Papa"
$ws.Cells.Item(332, 8).Value = "Pukará"
$ws.Cells.Item(332, 9).Value = "1a nueva(o)"
$ws.Cells.Item(332, 10).Value = 120
$ws.Cells.Item(332, 11).Value = 8500
$ws.Cells.Item(332, 12).Value = 9000
$ws.Cells.Item(332, 13).Value = 8750
$ws.Cells.Item(332, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(332, 15).Value = "Región de Ñuble"
$ws.Cells.Item(332, 16).Value = 350
$ws.Cells.Item(332, 17).Value = 25
$ws.Cells.Item(332, 18).Value = "Hortaliza"
